$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new header row at the very top; this shifts all existing
# data rows (and the two trailing blank rows) down by one, matching
# the diff (old row N -> new row N+1, dimension A1:B15 -> A1:B16).
$ws.Rows.Item(1).Insert()

# Populate the new header row with "Description" / "Value" labels
# (these become new shared-string entries, indices 12 and 13).
$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "Value"
